$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date value (row 8, column B) ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2023-03-27T14:49:55+00:00"

# --- Include ValueSets 2 sheet: DepartementOM-ROR URL (row 2, column A) ---
$wsDept = $wb.Worksheets.Item("Include ValueSets 2")
$wsDept.Range("A2").Value = "https://mos.esante.gouv.fr/NOS/JDV_J248-DepartementOM-ROR/FHIR/JDV-J248-DepartementOM-ROR"

# --- Include ValueSets 4 sheet: TerritoireSante-ROR URL (row 2, column A) ---
$wsTerr = $wb.Worksheets.Item("Include ValueSets 4")
$wsTerr.Range("A2").Value = "https://mos.esante.gouv.fr/NOS/JDV_J249-TerritoireSante-ROR/FHIR/JDV-J249-TerritoireSante-ROR"

# --- Include ValueSets 5 sheet: Pays-ROR URL (row 2, column A) ---
$wsPays = $wb.Worksheets.Item("Include ValueSets 5")
$wsPays.Range("A2").Value = "https://mos.esante.gouv.fr/NOS/JDV_J247-Pays-ROR/FHIR/JDV-J247-Pays-ROR"
